$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C3").Value = 153519
$ws.Range("C4").Value = 145060
$ws.Range("C5").Value = 8459
$ws.Range("C8").Value = 63.96
